$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents (values) of columns A and B for rows 1-10
for ($r = 1; $r -le 10; $r++) {
    $aVal = $ws.Cells.Item($r, 1).Value()
    $bVal = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r, 1).Value = $bVal
    $ws.Cells.Item($r, 2).Value = $aVal
}

# Swap the column widths of columns A and B (A was the wider column, now B is)
$ws.Columns.Item(1).ColumnWidth = 1.333333
$ws.Columns.Item(2).ColumnWidth = 2.333333
